# Update countries & provincias Spain
# Refresh of the COVID "Pais" data sheet:
#  - timestamp (A1) bumped from 07:32 to 08:49
#  - Uzbekistan overtakes Suiza (rows 60/61 swap rank, each with updated counts)
#  - Birmania overtakes Botsuana/Benin/Malta/Islandia (rows 145-149 shift rank)
#  - several other rows get refreshed case/death counters (rows 4,74,75,92,142,152)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- updated "last refreshed" timestamp -----------------------------------
$ws.Range("A1").Value = "Datos actualizados a 11 de Septiembre de 2020 a las 08:49"

# --- rows whose country name (column A) also changes rank ------------------
$ws.Range("A60").Value  = "Uzbekistan"
$ws.Range("A61").Value  = "Suiza"

$ws.Range("A145").Value = "Birmania"
$ws.Range("A146").Value = "Botsuana"
$ws.Range("A147").Value = "Benin"
$ws.Range("A148").Value = "Malta"
$ws.Range("A149").Value = "Islandia"

# --- refreshed numeric data (B=Casos totales, C=Nuevos casos, D=Casos activos,
#     E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes) -----------
function Set-RowData($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Range("B$row").Value = $b
    $ws.Range("C$row").Value = $c
    $ws.Range("D$row").Value = $d
    $ws.Range("E$row").Value = $e
    $ws.Range("F$row").Value = $f
    $ws.Range("G$row").Value = $g
    $ws.Range("H$row").Value = $h
}

Set-RowData 4   6588181 18  3880153 2511697 0 3 196331   # Estados Unidos
Set-RowData 60  45927   454 42555   2997    0 5 375      # Uzbekistan (now row 60)
Set-RowData 61  45711   0   38100   5591    0 0 2020     # Suiza (now row 61)
Set-RowData 74  26688   0   17106   8805    0 3 777      # El Salvador
Set-RowData 75  26564   51  23216   2551    0 9 797      # Australia
Set-RowData 92  11867   0   10371   1231    0 0 265      # Noruega
Set-RowData 142 2510    0   1313    1183    0 1 14       # Reunion
Set-RowData 145 2265    115 625     1626    0 0 14       # Birmania (now row 145)
Set-RowData 146 2252    0   546     1696    0 0 10       # Botsuana (now row 146)
Set-RowData 147 2242    0   1793    409     0 0 40       # Benin (now row 147)
Set-RowData 148 2204    0   1803    387     0 0 14       # Malta (now row 148)
Set-RowData 149 2157    0   2072    75      0 0 10       # Islandia (now row 149)
Set-RowData 152 1917    87  1354    544     0 0 19       # Georgia
